# Chain Lightning action and item created
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (19) for the "Chain Lightning" weapon item, copying the
# formatting of the existing "Fireball" weapon row (row 4) so that the new
# row inherits the same fills/styles used by other "Weapon" slot items.
$ws.Range("A4:G4").Copy($ws.Range("A19"))
$ws.Range("L4:M4").Copy($ws.Range("L19"))

# Fill in the new item's data.
$ws.Range("A19").Value = "Chain Lightning"
$ws.Range("B19").Value = "Weapon_ChainLightning"
$ws.Range("C19").Value = "Shoot a beam of lightning at a unit, that bounces onto nearby enemies."

# Update the active selection as it was left after the edit.
$ws.Range("R24").Select()
